$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.197.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.648.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.256'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.881.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.645.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.538'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.188.72'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").Value = '  -1.13%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  +3.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E30").Value = '  -1.28%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.271.79'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  +2.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.544'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.844'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.791.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("E47").Value = '  -1.06%  '
$ws.Range("E48").Value = '  +15.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0514'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.73%  '
